$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B22").Value = 6286
$ws.Range("C22").Value = 988
$ws.Range("D22").Value = 5759020
$ws.Range("E22").Value = 916.1660833598473
$ws.Range("F22").Value = 8.211396109485293
$ws.Range("G22").Value = 3.347280334728042
$ws.Range("H22").Value = 25.24073800224427
